# Release v2.3 deployment edit for dementia_Freiburg_v1 sheet:
# - Remove the conceptPath values that were filled in for PID_PSEUDONYMOUS (row 3)
#   and EXAMINATION_DATE (row 5); these two cells become blank again.
# - Update the AGE row's mapFunction / mapCDE text to add the new "subjectage"
#   (AGE*12) mapping alongside the existing ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").ClearContents()
$ws.Range("J5").ClearContents()

$ws.Range("L16").Value = "{stays the same}, {AGE*12}, {corresponds to one of the groups: {“-50y”},{”50-59y”},{”60-69y”},{”70-79y”},{”+80y”}}"
$ws.Range("M16").Value = "subjectageyears, subjectage, agegroup"

$ws.Range("J2").Select()
